# Finished the analytics part: refresh the sector-breakdown metrics
# (returns, cap, beta/alpha/sharpe/treynor, drawdown, std dev, R^2,
# expected return) with the latest computed figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells that become blank (numeric -> blank)
$ws.Cells.Item(7, 9).ClearContents()  # I7
$ws.Cells.Item(10, 9).ClearContents()  # I10
$ws.Cells.Item(11, 9).ClearContents()  # I11
$ws.Cells.Item(12, 9).ClearContents()  # I12

# Numeric value updates
# Row 2
$ws.Cells.Item(2, 2).Value = -13.74
$ws.Cells.Item(2, 3).Value = -5.307
$ws.Cells.Item(2, 4).Value = -19.793
$ws.Cells.Item(2, 5).Value = -2.772
$ws.Cells.Item(2, 6).Value = -10.347
$ws.Cells.Item(2, 7).Value = -7.474
$ws.Cells.Item(2, 8).Value = -10.866
$ws.Cells.Item(2, 9).Value = -7.366
$ws.Cells.Item(2, 10).Value = -1.97
$ws.Cells.Item(2, 11).Value = -11.217
$ws.Cells.Item(2, 12).Value = -6.293
# Row 3
$ws.Cells.Item(3, 2).Value = -15.363
$ws.Cells.Item(3, 3).Value = -5.636
$ws.Cells.Item(3, 4).Value = -21.774
$ws.Cells.Item(3, 5).Value = 5.745
$ws.Cells.Item(3, 6).Value = -9.612
$ws.Cells.Item(3, 7).Value = -4.655
$ws.Cells.Item(3, 8).Value = -6.598
$ws.Cells.Item(3, 9).Value = 2.533
$ws.Cells.Item(3, 10).Value = 0.091
$ws.Cells.Item(3, 11).Value = -3.924
$ws.Cells.Item(3, 12).Value = -2.214
# Row 4
$ws.Cells.Item(4, 2).Value = -13.807
$ws.Cells.Item(4, 3).Value = -5.562
$ws.Cells.Item(4, 4).Value = -24.224
$ws.Cells.Item(4, 5).Value = 0.825
$ws.Cells.Item(4, 6).Value = -9.202
$ws.Cells.Item(4, 7).Value = -5.167
$ws.Cells.Item(4, 8).Value = -7.841
$ws.Cells.Item(4, 9).Value = -0.593
$ws.Cells.Item(4, 10).Value = -2.782
$ws.Cells.Item(4, 11).Value = -7.328
$ws.Cells.Item(4, 12).Value = -2.96
# Row 5
$ws.Cells.Item(5, 2).Value = -13.424
$ws.Cells.Item(5, 3).Value = 10.607
$ws.Cells.Item(5, 4).Value = -29.568
$ws.Cells.Item(5, 5).Value = 22.34
$ws.Cells.Item(5, 6).Value = 2.433
$ws.Cells.Item(5, 7).Value = 5.527
$ws.Cells.Item(5, 8).Value = 14.033
$ws.Cells.Item(5, 9).Value = -76.22499999999999
$ws.Cells.Item(5, 10).Value = 6.107
$ws.Cells.Item(5, 11).Value = 13.438
$ws.Cells.Item(5, 12).Value = -11.549
# Row 6
$ws.Cells.Item(6, 2).Value = -17.161
$ws.Cells.Item(6, 3).Value = 34.593
$ws.Cells.Item(6, 4).Value = 14.897
$ws.Cells.Item(6, 5).Value = -7.182
$ws.Cells.Item(6, 6).Value = 40.312
$ws.Cells.Item(6, 7).Value = -9.246
$ws.Cells.Item(6, 8).Value = -17.719
$ws.Cells.Item(6, 9).Value = -499.751
$ws.Cells.Item(6, 10).Value = 15.991
$ws.Cells.Item(6, 11).Value = 57.695
$ws.Cells.Item(6, 12).Value = -10.639
# Row 7
$ws.Cells.Item(7, 2).Value = -4.178
$ws.Cells.Item(7, 3).Value = 6.967
$ws.Cells.Item(7, 4).Value = 3.198
$ws.Cells.Item(7, 5).Value = -1.675
$ws.Cells.Item(7, 6).Value = 7.981
$ws.Cells.Item(7, 7).Value = -2.175
$ws.Cells.Item(7, 8).Value = -4.325
$ws.Cells.Item(7, 10).Value = 3.42
$ws.Cells.Item(7, 11).Value = 10.878
$ws.Cells.Item(7, 12).Value = -2.518
# Row 8
$ws.Cells.Item(8, 2).Value = 64367.21
$ws.Cells.Item(8, 3).Value = 166301.59
$ws.Cells.Item(8, 4).Value = 66788.39
$ws.Cells.Item(8, 5).Value = 46738.05
$ws.Cells.Item(8, 6).Value = 198680.39
$ws.Cells.Item(8, 7).Value = 135461.86
$ws.Cells.Item(8, 8).Value = 144106.07
$ws.Cells.Item(8, 9).Value = 74203.16
$ws.Cells.Item(8, 10).Value = 48487.82
$ws.Cells.Item(8, 11).Value = 260273.99
$ws.Cells.Item(8, 12).Value = 22371.09
# Row 9
$ws.Cells.Item(9, 2).Value = 0.722
$ws.Cells.Item(9, 3).Value = 0.996
$ws.Cells.Item(9, 4).Value = 0.9389999999999999
$ws.Cells.Item(9, 5).Value = 0.619
$ws.Cells.Item(9, 6).Value = 0.917
$ws.Cells.Item(9, 7).Value = 0.996
$ws.Cells.Item(9, 8).Value = 0.986
$ws.Cells.Item(9, 9).Value = -0.028
$ws.Cells.Item(9, 10).Value = 0.663
$ws.Cells.Item(9, 11).Value = 1.101
$ws.Cells.Item(9, 12).Value = 0.468
# Row 10
$ws.Cells.Item(10, 2).Value = -10.318
$ws.Cells.Item(10, 3).Value = -1.008
$ws.Cells.Item(10, 4).Value = -4.394
$ws.Cells.Item(10, 5).Value = -7.122
$ws.Cells.Item(10, 6).Value = 0.539
$ws.Cells.Item(10, 7).Value = -10.152
$ws.Cells.Item(10, 8).Value = -12.231
$ws.Cells.Item(10, 10).Value = -2.324
$ws.Cells.Item(10, 11).Value = 2.203
$ws.Cells.Item(10, 12).Value = -6.953
# Row 11
$ws.Cells.Item(11, 2).Value = -31.689
$ws.Cells.Item(11, 3).Value = 34.964
$ws.Cells.Item(11, 4).Value = 3.146
$ws.Cells.Item(11, 5).Value = -11.166
$ws.Cells.Item(11, 6).Value = 45.694
$ws.Cells.Item(11, 7).Value = -16.241
$ws.Cells.Item(11, 8).Value = -25.072
$ws.Cells.Item(11, 10).Value = 17.888
$ws.Cells.Item(11, 11).Value = 49.357
$ws.Cells.Item(11, 12).Value = -23.774
# Row 12
$ws.Cells.Item(12, 2).Value = -0.076
$ws.Cells.Item(12, 3).Value = 0.057
$ws.Cells.Item(12, 4).Value = 0.02
$ws.Cells.Item(12, 5).Value = -0.048
$ws.Cells.Item(12, 6).Value = 0.073
$ws.Cells.Item(12, 7).Value = -0.035
$ws.Cells.Item(12, 8).Value = -0.057
$ws.Cells.Item(12, 10).Value = 0.032
$ws.Cells.Item(12, 11).Value = 0.08699999999999999
$ws.Cells.Item(12, 12).Value = -0.082
# Row 13
$ws.Cells.Item(13, 2).Value = -35.044
$ws.Cells.Item(13, 3).Value = -25.478
$ws.Cells.Item(13, 4).Value = -36.952
$ws.Cells.Item(13, 5).Value = -48.448
$ws.Cells.Item(13, 6).Value = -21.505
$ws.Cells.Item(13, 7).Value = -26.76
$ws.Cells.Item(13, 8).Value = -43.575
$ws.Cells.Item(13, 9).Value = -1369.897
$ws.Cells.Item(13, 10).Value = -23.868
$ws.Cells.Item(13, 11).Value = -33.072
$ws.Cells.Item(13, 12).Value = -22.812
# Row 14
$ws.Cells.Item(14, 2).Value = 17.296
$ws.Cells.Item(14, 3).Value = 16.216
$ws.Cells.Item(14, 4).Value = 60.363
$ws.Cells.Item(14, 5).Value = 26.659
$ws.Cells.Item(14, 6).Value = 14.628
$ws.Cells.Item(14, 7).Value = 21.41
$ws.Cells.Item(14, 8).Value = 22.445
$ws.Cells.Item(14, 9).Value = 631.391
$ws.Cells.Item(14, 10).Value = 11.861
$ws.Cells.Item(14, 11).Value = 19.415
$ws.Cells.Item(14, 12).Value = 16.082
# Row 15
$ws.Cells.Item(15, 2).Value = 0.201
$ws.Cells.Item(15, 3).Value = 0.622
$ws.Cells.Item(15, 4).Value = 0.081
$ws.Cells.Item(15, 5).Value = 0.067
$ws.Cells.Item(15, 6).Value = 0.666
$ws.Cells.Item(15, 7).Value = 0.32
$ws.Cells.Item(15, 8).Value = 0.252
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 0.672
$ws.Cells.Item(15, 11).Value = 0.663
$ws.Cells.Item(15, 12).Value = 0.267
# Row 16
$ws.Cells.Item(16, 2).Value = 6.14
$ws.Cells.Item(16, 3).Value = 7.975
$ws.Cells.Item(16, 4).Value = 7.592
$ws.Cells.Item(16, 5).Value = 5.447
$ws.Cells.Item(16, 6).Value = 7.442
$ws.Cells.Item(16, 7).Value = 7.976
$ws.Cells.Item(16, 8).Value = 7.906
$ws.Cells.Item(16, 9).Value = 1.114
$ws.Cells.Item(16, 10).Value = 8.676
$ws.Cells.Item(16, 11).Value = 8.676
$ws.Cells.Item(16, 12).Value = 4.435
